# Update cryptos list data (prices + volume deltas) to match
# "Updated cryptos list on Fri May 24 08:50:33 UTC 2024 with GitHub Actions"
#
# Numeric-looking price strings are written with a leading apostrophe
# (quote-prefix) so Excel keeps them as text, preserving formatting such
# as trailing zeros (e.g. "1.00", "12.10") instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.082.88'
$ws.Range("E2").Value = '  -3.68%  '
$ws.Range("D3").Value = '3.671.99'
$ws.Range("E3").Value = '  -3.41%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''593.09'
$ws.Range("E5").Value = '  -3.33%  '
$ws.Range("D6").Value = '''164.94'
$ws.Range("E6").Value = '  -6.98%  '
$ws.Range("D7").Value = '3.669.66'
$ws.Range("E7").Value = '  -3.46%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''0.523'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("E10").Value = '  -5.14%  '
$ws.Range("D11").Value = '''6.13'
$ws.Range("E11").Value = '  -5.54%  '
$ws.Range("D12").Value = '''0.460'
$ws.Range("E12").Value = '  -4.94%  '
$ws.Range("D13").Value = '''37.36'
$ws.Range("E13").Value = '  -6.03%  '
$ws.Range("E14").Value = '  -6.26%  '
$ws.Range("D15").Value = '4.283.70'
$ws.Range("E15").Value = '  -3.36%  '
$ws.Range("D16").Value = '3.666.27'
$ws.Range("E16").Value = '  -3.45%  '
$ws.Range("D17").Value = '67.069.49'
$ws.Range("E17").Value = '  -3.74%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '''7.13'
$ws.Range("E18").Value = '  -5.61%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '''0.114'
$ws.Range("E19").Value = '  -4.21%  '
$ws.Range("D20").Value = '''17.06'
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("D21").Value = '''488.57'
$ws.Range("E21").Value = '  -3.55%  '
$ws.Range("D22").Value = '''9.05'
$ws.Range("E22").Value = '  -5.88%  '
$ws.Range("D23").Value = '''0.714'
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("D24").Value = '''85.23'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("E25").Value = '  -7.61%  '
$ws.Range("E26").Value = '  -5.40%  '
$ws.Range("D27").Value = '''12.10'
$ws.Range("E27").Value = '  -4.15%  '
$ws.Range("D28").Value = '''0.996'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("D29").Value = '''9.91'
$ws.Range("E29").Value = '  -6.14%  '
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("E31").Value = '  -6.74%  '
$ws.Range("D32").Value = '''7.69'
$ws.Range("E32").Value = '  -4.19%  '
$ws.Range("D33").Value = '''31.55'
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = '3.807.59'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("D35").Value = '3.608.21'
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("E36").Value = '  -6.68%  '
$ws.Range("D37").Value = '''0.998'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '''0.988'
$ws.Range("E38").Value = '  -5.38%  '
$ws.Range("D39").Value = '''5.73'
$ws.Range("E39").Value = '  -6.25%  '
$ws.Range("D40").Value = '''0.131'
$ws.Range("E40").Value = '  -7.48%  '
$ws.Range("E41").Value = '  -5.07%  '
$ws.Range("D42").Value = '''437.09'
$ws.Range("E42").Value = '  -9.12%  '
$ws.Range("D43").Value = '''48.57'
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("D44").Value = '''1.92'
$ws.Range("E44").Value = '  -6.94%  '
$ws.Range("E45").Value = '  -9.30%  '
$ws.Range("D46").Value = '''8.30'
$ws.Range("E46").Value = '  -3.11%  '
$ws.Range("D48").Value = '''142.06'
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("D49").Value = '''39.65'
$ws.Range("E49").Value = '  -10.07%  '
$ws.Range("D50").Value = '2.745.30'
$ws.Range("E50").Value = '  -6.29%  '
$ws.Range("D51").Value = '''0.0345'
$ws.Range("E51").Value = '  -4.79%  '
